$d = $word.ActiveDocument

# Target paragraph is the last (empty) paragraph, styled "ListParagraph" with no text.
$headingPara = $d.Paragraphs.Last
$headingPara.Range.Text = "Building scrolling experiences in Flutter:"
$headingPara.Style = "Heading 1"

$clearRange = $d.Range($headingPara.Range.Start, $headingPara.Range.End)
$clearRange.Text = ""
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Building scrolling experiences in Flutter</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insPoint = $d.Range($headingPara.Range.Start, $headingPara.Range.Start)
$insPoint.InsertXML($xml)

$r = $headingPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$firstItem = $d.Paragraphs.Last
$firstItem.Range.Text = "CustomScrollView (in body of scaffold)"
$firstItem.Style = "List Paragraph"
$gal = $word.ListGalleries.Item(1)
$lt = $gal.ListTemplates.Item(1)
$firstItem.Range.ListFormat.ApplyListTemplate($lt)

$clearRange = $d.Range($firstItem.Range.Start, $firstItem.Range.End)
$clearRange.Text = ""
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>CustomScrollView</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (in body of scaffold)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insPoint = $d.Range($firstItem.Range.Start, $firstItem.Range.Start)
$insPoint.InsertXML($xml)

# ---- item 2 ----
$prevItem = $d.Paragraphs.Last
$pr = $prevItem.Range
$pr.Collapse(0)
$pr.InsertParagraphAfter()
$curItem = $d.Paragraphs.Last
$curItem.Range.Text = "SliverAppbar (pinned, stretched, flexible space etc. )"
$clearRange = $d.Range($curItem.Range.Start, $curItem.Range.End)
$clearRange.Text = ""
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>SliverAppbar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (pinned, stretched, flexible space </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>etc. )</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insPoint = $d.Range($curItem.Range.Start, $curItem.Range.Start)
$insPoint.InsertXML($xml)

# ---- item 3 ----
$prevItem = $d.Paragraphs.Last
$pr = $prevItem.Range
$pr.Collapse(0)
$pr.InsertParagraphAfter()
$curItem = $d.Paragraphs.Last
$curItem.Range.Text = "Flexible spacebar (collapse modes, strectch modes)"
$clearRange = $d.Range($curItem.Range.Start, $curItem.Range.End)
$clearRange.Text = ""
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Flexible spacebar (collapse modes, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>strectch</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> modes)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insPoint = $d.Range($curItem.Range.Start, $curItem.Range.Start)
$insPoint.InsertXML($xml)

# ---- item 4 ----
$prevItem = $d.Paragraphs.Last
$pr = $prevItem.Range
$pr.Collapse(0)
$pr.InsertParagraphAfter()
$curItem = $d.Paragraphs.Last
$curItem.Range.Text = "Box Decoration (gradients)"
$clearRange = $d.Range($curItem.Range.Start, $curItem.Range.End)
$clearRange.Text = ""
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Box Decoration (gradients)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insPoint = $d.Range($curItem.Range.Start, $curItem.Range.Start)
$insPoint.InsertXML($xml)

# ---- item 5 ----
$prevItem = $d.Paragraphs.Last
$pr = $prevItem.Range
$pr.Collapse(0)
$pr.InsertParagraphAfter()
$curItem = $d.Paragraphs.Last
$curItem.Range.Text = "Sliver List (delegate : SliverChildBuilderDelegate)"
$clearRange = $d.Range($curItem.Range.Start, $curItem.Range.End)
$clearRange.Text = ""
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Sliver List (</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>delegate :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>SliverChildBuilderDelegate</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insPoint = $d.Range($curItem.Range.Start, $curItem.Range.Start)
$insPoint.InsertXML($xml)

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    Write-Output ("[$i] " + $pp.Range.Text)
}
